$d = $word.ActiveDocument

# The page is being trimmed of its site-footer boilerplate:
#   - the blank paragraph right before "Ver no Jupiter ..."
#   - the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph
#   - the "© 2020 . Contact: ..." paragraph
# These three paragraphs are removed as a block, leaving the remaining
# blank paragraph (and the page-break paragraph after it) untouched.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter*Salvar em pdf*Salvar em docx*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $prev = $target.Previous()
    $next = $target.Next()

    $startPos = $target.Range.Start
    $endPos = $target.Range.End

    if ($prev -ne $null -and $prev.Range.Text.Trim() -eq "") {
        $startPos = $prev.Range.Start
    }

    if ($next -ne $null -and $next.Range.Text -like "*Contact*") {
        $endPos = $next.Range.End
    }

    $deleteRange = $d.Range($startPos, $endPos)
    $deleteRange.Delete()
}
